# Migrating to use yaml data format
#
# Adds four new worksheets - DACCS_pars, DACCS_sets, DACCS_region, DACCS_time -
# carrying the DACCS_LT technology parameters/sets/region-factors/time-factors
# that used to live only inline in Sheet1.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$pars   = $wb.Worksheets.Add($null, $sheet1)
$pars.Name = "DACCS_pars"
$sets   = $wb.Worksheets.Add($null, $pars)
$sets.Name = "DACCS_sets"
$region = $wb.Worksheets.Add($null, $sets)
$region.Name = "DACCS_region"
$time   = $wb.Worksheets.Add($null, $region)
$time.Name = "DACCS_time"

# ---------------------------------------------------------------------------
# DACCS_pars : name | par | unit | value
# ---------------------------------------------------------------------------
$pars.Range("B1").Value = "par"
$pars.Range("D1").Value = "value"
$pars.Range("A1").Value = "name"

$pars.Range("B2").Value = "inv_cost"
$pars.Range("D2").Value = 2500
$pars.Range("A2").Value = "DACCS_LT"

$pars.Range("B3").Value = "var_cost"
$pars.Range("D3").Value = 100
$pars.Range("A3").Value = "DACCS_LT"

$pars.Range("B4").Value = "fix_cost"
$pars.Range("D4").Value = 10
$pars.Range("A4").Value = "DACCS_LT"

$pars.Range("B5").Value = "input"
$pars.Range("D5").Value = 1
$pars.Range("A5").Value = "DACCS_LT"

$pars.Range("B6").Value = "output"
$pars.Range("D6").Value = 0
$pars.Range("A6").Value = "DACCS_LT"

$pars.Range("B7").Value = "emission_factor"
$pars.Range("D7").Value = -20
$pars.Range("A7").Value = "DACCS_LT"

$pars.Range("B8").Value = "capacity_factor"
$pars.Range("D8").Value = 1
$pars.Range("A8").Value = "DACCS_LT"

$pars.Range("B9").Value = "technical_lifetime"
$pars.Range("D9").Value = 20
$pars.Range("A9").Value = "DACCS_LT"

$pars.Range("B10").Value = "initial_new_capacity_up"
$pars.Range("D10").Value = 0.5
$pars.Range("A10").Value = "DACCS_LT"

$pars.Range("B11").Value = "growth_new_capacity_up"
$pars.Range("D11").Value = 0.5
$pars.Range("A11").Value = "DACCS_LT"

# ---------------------------------------------------------------------------
# DACCS_sets : name | set | value
# ---------------------------------------------------------------------------
$sets.Range("A1").Value = "name"
$sets.Range("B1").Value = "set"
$sets.Range("C1").Value = "value"

$sets.Range("A2").Value = "DACCS_LT"
$sets.Range("B2").Value = "time"
$sets.Range("C2").Value = "year"

$sets.Range("A3").Value = "DACCS_LT"
$sets.Range("B3").Value = "mode"
$sets.Range("C3").Value = "standard"

$sets.Range("A4").Value = "DACCS_LT"
$sets.Range("B4").Value = "emission"
$sets.Range("C4").Value = "CO2"

$sets.Range("A5").Value = "DACCS_LT"
$sets.Range("B5").Value = "commodity_in"
$sets.Range("C5").Value = "electricity"

$sets.Range("A6").Value = "DACCS_LT"
$sets.Range("B6").Value = "level_in"
$sets.Range("C6").Value = "final"

$sets.Range("A7").Value = "DACCS_LT"
$sets.Range("B7").Value = "commodity_out"
$sets.Range("C7").Value = "light"

$sets.Range("A8").Value = "DACCS_LT"
$sets.Range("B8").Value = "level_out"
$sets.Range("C8").Value = "useful"

$sets.Range("A9").Value = "DACCS_LT"
$sets.Range("B9").Value = "first_active_year"
$sets.Range("C9").Value = 2020

# ---------------------------------------------------------------------------
# DACCS_region : name | region | unit | inv_cost
# ---------------------------------------------------------------------------
$region.Range("A1").Value = "name"
$region.Range("B1").Value = "region"
$region.Range("D1").Value = "inv_cost"

$region.Range("A2").Value = "DACCS_LT"
$region.Range("B2").Value = "NAM"
$region.Range("D2").Value = 1

$region.Range("A3").Value = "DACCS_LT"
$region.Range("B3").Value = "WEU"
$region.Range("D3").Value = 1.2

# ---------------------------------------------------------------------------
# DACCS_time : name | unit | inv_cost_reduction_rate
# ---------------------------------------------------------------------------
$time.Range("A1").Value = "name"
$time.Range("C1").Value = "inv_cost_reduction_rate"

$time.Range("A2").Value = "DACCS_LT"
$time.Range("C2").Value = 0.05

# ---------------------------------------------------------------------------
# "unit" annotation columns, filled in afterwards
# ---------------------------------------------------------------------------
$pars.Range("C1").Value = "unit"
$pars.Range("C2").Value = "$/kW"
$pars.Range("C3").Value = "$/kWa"
$pars.Range("C4").Value = "$/kWa"
$pars.Range("C7").Value = "tCO2/kWa"
$pars.Range("C9").Value = "y"
$pars.Range("C10").Value = "GW"

$time.Range("B1").Value = "unit"
$time.Range("B2").Value = "% / yr"

$region.Range("C1").Value = "unit"
$region.Range("C2").Value = "fraction"
$region.Range("C3").Value = "fraction"

# ---------------------------------------------------------------------------
# Column widths (bestFit, as saved by Excel for the new sheets)
# ---------------------------------------------------------------------------
$pars.Columns.Item(1).ColumnWidth = 8.5
$pars.Columns.Item(2).ColumnWidth = 21.5

$sets.Columns.Item(1).ColumnWidth = 8.5
$sets.Columns.Item(2).ColumnWidth = 13.35

# ---------------------------------------------------------------------------
# Selections / active sheet + tab bookkeeping
# ---------------------------------------------------------------------------
$sheet1.Range("A29").Select()

$sets.Range("B7").Select()
$region.Range("C4").Select()
$time.Range("B3").Select()

$pars.Range("G17").Select()
$pars.Activate()
